$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: placeholder fill so every target row/cell exists with a plain (non-numeric) value ---
# (actual values are written in step 4; writing the real "91.88%"-style strings before the
#  column's NumberFormat is forced to Text would make Excel auto-convert them to percentage numbers)
$ws.Range("A2:K7").Value = "_"

# --- Step 2: normalize the whole new block to the same style used by the rest of the table ---
$ws.Range("A1").Copy()
$ws.Range("A2:K7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: force Text interpretation for the numeric-looking percentage cells ---
$ws.Range("B2:K7").NumberFormat = "@"

# --- Step 4: write the real row labels + result values ---
$ws.Range("A2").Value = "pre_att_all"
$ws.Range("B2").Value = "91.88%"
$ws.Range("C2").Value = "91.00%"
$ws.Range("D2").Value = "86.62%"
$ws.Range("E2").Value = "64.70%"
$ws.Range("F2").Value = "54.00%"
$ws.Range("G2").Value = "51.21%"
$ws.Range("H2").Value = "50.62%"
$ws.Range("I2").Value = "50.16%"
$ws.Range("J2").Value = "50.06%"
$ws.Range("K2").Value = "50.07%"
$ws.Range("A3").Value = "pre_att_cls"
$ws.Range("B3").Value = "92.10%"
$ws.Range("C3").Value = "92.10%"
$ws.Range("D3").Value = "91.96%"
$ws.Range("E3").Value = "91.76%"
$ws.Range("F3").Value = "90.88%"
$ws.Range("G3").Value = "89.49%"
$ws.Range("H3").Value = "88.06%"
$ws.Range("I3").Value = "86.26%"
$ws.Range("J3").Value = "84.26%"
$ws.Range("K3").Value = "81.83%"
$ws.Range("A4").Value = "post_att_all"
$ws.Range("B4").Value = "92.05%"
$ws.Range("C4").Value = "91.73%"
$ws.Range("D4").Value = "90.81%"
$ws.Range("E4").Value = "86.95%"
$ws.Range("F4").Value = "74.22%"
$ws.Range("G4").Value = "58.20%"
$ws.Range("H4").Value = "52.30%"
$ws.Range("I4").Value = "50.76%"
$ws.Range("J4").Value = "50.06%"
$ws.Range("K4").Value = "50.02%"
$ws.Range("A5").Value = "post_att_cls"
$ws.Range("B5").Value = "92.11%"
$ws.Range("C5").Value = "92.11%"
$ws.Range("D5").Value = "92.05%"
$ws.Range("E5").Value = "92.03%"
$ws.Range("F5").Value = "91.93%"
$ws.Range("G5").Value = "91.87%"
$ws.Range("H5").Value = "91.17%"
$ws.Range("I5").Value = "89.64%"
$ws.Range("J5").Value = "87.49%"
$ws.Range("K5").Value = "84.96%"
$ws.Range("A6").Value = "last_cls"
$ws.Range("B6").Value = "92.14%"
$ws.Range("C6").Value = "92.15%"
$ws.Range("D6").Value = "92.15%"
$ws.Range("E6").Value = "92.16%"
$ws.Range("F6").Value = "92.15%"
$ws.Range("G6").Value = "92.16%"
$ws.Range("H6").Value = "92.12%"
$ws.Range("I6").Value = "92.16%"
$ws.Range("J6").Value = "92.13%"
$ws.Range("K6").Value = "92.13%"
$ws.Range("A7").Value = "logits"
$ws.Range("B7").Value = "92.15%"
$ws.Range("C7").Value = "92.14%"
$ws.Range("D7").Value = "92.14%"
$ws.Range("E7").Value = "92.13%"
$ws.Range("F7").Value = "92.13%"
$ws.Range("G7").Value = "92.15%"
$ws.Range("H7").Value = "92.11%"
$ws.Range("I7").Value = "92.09%"
$ws.Range("J7").Value = "92.11%"
$ws.Range("K7").Value = "92.10%"
